# Add a new, empty "Title and Content" slide as slide #2 (right after the
# existing title slide). Layout 2 corresponds to slideLayout2.xml ("Title
# and Content"), matching the title + content placeholders added in the
# target slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Add(2, 2)
